$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round existing values for rows 145-149 to match updated source data
$ws.Range("B145").Value = 133
$ws.Range("C145").Value = 6.2

$ws.Range("B146").Value = 146.5
$ws.Range("C146").Value = 6.5

$ws.Range("B147").Value = 155.2
$ws.Range("C147").Value = 7.1

$ws.Range("B148").Value = 158.4
$ws.Range("C148").Value = 6.9

$ws.Range("B149").Value = 150.1
$ws.Range("C149").Value = 6.8

# Add new row 150 with the latest period data
$ws.Range("A150").Value = "06_21/22"
$ws.Range("B150").Value = 173.4
$ws.Range("C150").Value = 7.100000000000001
